# LevelTable.xlsx — add a new "Lv_0" level row (row 10) to Sheet1,
# mirroring the existing data rows (5-9) for values and formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row: ID=6, LvPfbName="Lv_0", InitMoney=315, MoneyColorId=3,
# LvDisplayName="Level.DisplayName6"
$ws.Range("B10").Value = 6
$ws.Range("D10").Value = "Lv_0"
$ws.Range("E10").Value = 315
$ws.Range("F10").Value = 3
$ws.Range("G10").Value = "Level.DisplayName6"

# Carry over the same cell formatting used by the rows above (D has the
# "Lv_x" style, G has the "Level.DisplayNameN" style).
$ws.Range("D9").Copy()
$ws.Range("D10").PasteSpecial(-4122)

$ws.Range("G9").Copy()
$ws.Range("G10").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Keep the active-cell selection near the new data, matching the saved
# session state after the edit.
$ws.Range("H15").Select() | Out-Null
